$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 9 ("Are we ready?") so the new
# problem gets its own row while the rows below shift down by one.
$ws.Rows("9:9").Insert()

# New problem entry.
$ws.Range("A9").Value = "8. Reviews for 2nd slide"

# Match the saved selection state from the edit.
$ws.Range("A11").Select()
